$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 27, pushing existing rows 27-56 down to 28-57.
$ws.Rows.Item(27).Insert()

# Populate the newly inserted row 27 with the new record.
$ws.Range("A27").Value = 1
$ws.Range("B27").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C27").Value = "Arica y Parinacota"
$ws.Range("D27").Value = 44622
$ws.Range("E27").Value = 15
$ws.Range("F27").Value = 100112009
$ws.Range("G27").Value = "Acelga"
$ws.Range("H27").Value = "Sin especificar"
$ws.Range("I27").Value = "Segunda"
$ws.Range("J27").Value = 200
$ws.Range("K27").Value = 900
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = 950
$ws.Range("N27").Value = "$/atado 2,5 a 3 kilos"
$ws.Range("O27").Value = "Región de Arica y Parinacota"
$ws.Range("P27").Value = 317
$ws.Range("Q27").Value = 3
$ws.Range("R27").Value = "Hortaliza"
